$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("PPRiFUfIIaIoE")

$cell = $ws.Range("B1")
$cell.Value = "Pot Perc Red in Fuel Use (dimensionless)"
$cell.WrapText = $true

$ws.Activate()
$ws.Range("B1").Select() | Out-Null

$wb.Worksheets.Item("About").Activate() | Out-Null
